$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on the data range so numeric-looking strings
# (e.g. "577.04") are not auto-converted to numbers by Excel,
# matching the original inline-string cell type.
$dataRange = $ws.Range("D2:E51")
$dataRange.NumberFormat = "@"

# Row 2
$ws.Range("D2").Value = '70.648.27'
$ws.Range("E2").Value = '  +0.71%  '

# Row 3
$ws.Range("D3").Value = '3.645.57'
$ws.Range("E3").Value = '  +5.98%  '

# Row 4
$ws.Range("E4").Value = '  +0.09%  '

# Row 5
$ws.Range("D5").Value = '577.04'
$ws.Range("E5").Value = '  -1.28%  '

# Row 6
$ws.Range("D6").Value = '176.00'
$ws.Range("E6").Value = '  -1.24%  '

# Row 7
$ws.Range("D7").Value = '3.635.13'
$ws.Range("E7").Value = '  +5.90%  '

# Row 8
$ws.Range("D8").Value = '0.611'
$ws.Range("E8").Value = '  +1.76%  '

# Row 9
$ws.Range("E9").Value = '  +0.17%  '

# Row 10
$ws.Range("E10").Value = '  -4.01%  '

# Row 11
$ws.Range("D11").Value = '6.83'
$ws.Range("E11").Value = '  +23.70%  '

# Row 12
$ws.Range("D12").Value = '0.602'
$ws.Range("E12").Value = '  +2.43%  '

# Row 13
$ws.Range("D13").Value = '48.47'
$ws.Range("E13").Value = '  -0.94%  '

# Row 14
$ws.Range("D14").Value = '0.0000287'
$ws.Range("E14").Value = '  +0.41%  '

# Row 15
$ws.Range("D15").Value = '4.232.54'
$ws.Range("E15").Value = '  +6.09%  '

# Row 16
$ws.Range("D16").Value = '669.08'
$ws.Range("E16").Value = '  -3.42%  '

# Row 17
$ws.Range("D17").Value = '8.86'
$ws.Range("E17").Value = '  +1.97%  '

# Row 18
$ws.Range("D18").Value = '3.641.34'
$ws.Range("E18").Value = '  +6.08%  '

# Row 19
$ws.Range("D19").Value = '70.718.78'
$ws.Range("E19").Value = '  +0.88%  '

# Row 20
$ws.Range("E20").Value = '  +0.52%  '

# Row 21
$ws.Range("E21").Value = '  +0.09%  '

# Row 22
$ws.Range("E22").Value = '  -0.78%  '

# Row 23
$ws.Range("E23").Value = '  +2.86%  '

# Row 24
$ws.Range("D24").Value = '17.10'
$ws.Range("E24").Value = '  -0.08%  '

# Row 25
$ws.Range("D25").Value = '100.38'
$ws.Range("E25").Value = '  -0.91%  '

# Row 26
$ws.Range("E26").Value = '  -1.11%  '

# Row 27
$ws.Range("E27").Value = '  +3.75%  '

# Row 28
$ws.Range("E28").Value = '  -0.10%  '

# Row 29
$ws.Range("D29").Value = '9.98'
$ws.Range("E29").Value = '  +3.64%  '

# Row 30
$ws.Range("D30").Value = '34.76'
$ws.Range("E30").Value = '  +3.39%  '

# Row 31
$ws.Range("E31").Value = '  +0.68%  '

# Row 32
$ws.Range("D32").Value = '8.99'
$ws.Range("E32").Value = '  +2.26%  '

# Row 33
$ws.Range("D33").Value = '1.39'
$ws.Range("E33").Value = '  -4.10%  '

# Row 34
$ws.Range("D34").Value = '7.29'
$ws.Range("E34").Value = '  +1.34%  '

# Row 35
$ws.Range("D35").Value = '3.99'
$ws.Range("E35").Value = '  +3.19%  '

# Row 36
$ws.Range("D36").Value = '582.92'
$ws.Range("E36").Value = '  +1.43%  '

# Row 37
$ws.Range("D37").Value = '11.04'
$ws.Range("E37").Value = '  -0.18%  '

# Row 38
$ws.Range("E38").Value = '  +3.10%  '

# Row 39
$ws.Range("D39").Value = '58.23'
$ws.Range("E39").Value = '  -0.89%  '

# Row 40
$ws.Range("D40").Value = '0.999'
$ws.Range("E40").Value = '  -0.09%  '

# Row 41
$ws.Range("D41").Value = '3.581.93'
$ws.Range("E41").Value = '  -0.15%  '

# Row 42
$ws.Range("D42").Value = '0.0453'
$ws.Range("E42").Value = '  +7.06%  '

# Row 43
$ws.Range("E43").Value = '  +0.98%  '

# Row 44
$ws.Range("E44").Value = '  +2.36%  '

# Row 45
$ws.Range("B45").Value = 'InjectiveProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D45").Value = '34.72'
$ws.Range("E45").Value = '  -1.78%  '

# Row 46
$ws.Range("B46").Value = 'PEPE'
$ws.Range("C46").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D46").Value = '0.0₃0743'
$ws.Range("E46").Value = '  +0.05%  '

# Row 47
$ws.Range("E47").Value = '  +0.36%  '

# Row 48
$ws.Range("D48").Value = '2.92'
$ws.Range("E48").Value = '  +8.59%  '

# Row 49
$ws.Range("D49").Value = '0.132'
$ws.Range("E49").Value = '  +2.21%  '

# Row 50
$ws.Range("D50").Value = '134.98'
$ws.Range("E50").Value = '  +1.03%  '

# Row 51
$ws.Range("D51").Value = '2.94'
$ws.Range("E51").Value = '  +7.32%  '

# Restore default (General) style on the data range so the written
# cells keep their text VALUES but do not retain the temporary
# text number-format / style index that was needed to stop Excel
# from re-typing the numeric-looking strings as numbers.
$dataRange.ClearFormats()
